# DOMA-2542 Localization for Excel template (contacts)
# Replace hard-coded Russian strings (sheet name + header row) with
# i18n placeholder tokens so the template can be rendered for any locale.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet name: "Контакты" -> "{d.i18n.sheetName}"
$ws.Name = "{d.i18n.sheetName}"

# Header row (row 1): localize each column header
$ws.Range("A1").Value = "{d.i18n.name}"      # was "Имя"
$ws.Range("B1").Value = "{d.i18n.address}"   # was "Адрес"
$ws.Range("C1").Value = "{d.i18n.unitName}"  # was "Квартира"
$ws.Range("D1").Value = "{d.i18n.phone}"     # was "Телефон"
$ws.Range("E1").Value = "{d.i18n.email}"     # was "Почта"
